$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.Value = "'59.323.42"
$cell.Style = "Normal"
$cell = $ws.Range("E2")
$cell.Value = "'  +0.47%  "
$cell.Style = "Normal"
$cell = $ws.Range("D3")
$cell.Value = "'2.606.41"
$cell.Style = "Normal"
$cell = $ws.Range("E3")
$cell.Value = "'  +0.45%  "
$cell.Style = "Normal"
$cell = $ws.Range("E4")
$cell.Value = "'  -0.01%  "
$cell.Style = "Normal"
$cell = $ws.Range("D5")
$cell.Value = "'536.85"
$cell.Style = "Normal"
$cell = $ws.Range("E5")
$cell.Value = "'  +3.46%  "
$cell.Style = "Normal"
$cell = $ws.Range("D6")
$cell.Value = "'140.44"
$cell.Style = "Normal"
$cell = $ws.Range("E6")
$cell.Value = "'  +0.90%  "
$cell.Style = "Normal"
$cell = $ws.Range("E7")
$cell.Value = "'  +0.32%  "
$cell.Style = "Normal"
$cell = $ws.Range("D8")
$cell.Value = "'0.567"
$cell.Style = "Normal"
$cell = $ws.Range("E8")
$cell.Value = "'  +0.06%  "
$cell.Style = "Normal"
$cell = $ws.Range("D9")
$cell.Value = "'2.613.92"
$cell.Style = "Normal"
$cell = $ws.Range("E9")
$cell.Value = "'  -0.25%  "
$cell.Style = "Normal"
$cell = $ws.Range("D10")
$cell.Value = "'6.44"
$cell.Style = "Normal"
$cell = $ws.Range("E10")
$cell.Value = "'  -0.65%  "
$cell.Style = "Normal"
$cell = $ws.Range("E11")
$cell.Value = "'  +1.81%  "
$cell.Style = "Normal"
$cell = $ws.Range("E13")
$cell.Value = "'  +1.76%  "
$cell.Style = "Normal"
$cell = $ws.Range("D14")
$cell.Value = "'3.067.75"
$cell.Style = "Normal"
$cell = $ws.Range("E14")
$cell.Value = "'  +0.27%  "
$cell.Style = "Normal"
$cell = $ws.Range("D15")
$cell.Value = "'59.250.48"
$cell.Style = "Normal"
$cell = $ws.Range("E15")
$cell.Value = "'  +0.40%  "
$cell.Style = "Normal"
$cell = $ws.Range("D16")
$cell.Value = "'20.53"
$cell.Style = "Normal"
$cell = $ws.Range("E16")
$cell.Value = "'  +0.80%  "
$cell.Style = "Normal"
$cell = $ws.Range("D17")
$cell.Value = "'2.607.93"
$cell.Style = "Normal"
$cell = $ws.Range("E17")
$cell.Value = "'  -0.27%  "
$cell.Style = "Normal"
$cell = $ws.Range("E18")
$cell.Value = "'  +0.83%  "
$cell.Style = "Normal"
$cell = $ws.Range("D19")
$cell.Value = "'342.47"
$cell.Style = "Normal"
$cell = $ws.Range("E19")
$cell.Value = "'  +0.96%  "
$cell.Style = "Normal"
$cell = $ws.Range("E20")
$cell.Value = "'  +0.91%  "
$cell.Style = "Normal"
$cell = $ws.Range("D21")
$cell.Value = "'10.09"
$cell.Style = "Normal"
$cell = $ws.Range("E21")
$cell.Value = "'  -0.72%  "
$cell.Style = "Normal"
$cell = $ws.Range("E22")
$cell.Value = "'  -1.29%  "
$cell.Style = "Normal"
$cell = $ws.Range("E23")
$cell.Value = "'  +0.06%  "
$cell.Style = "Normal"
$cell = $ws.Range("D24")
$cell.Value = "'67.48"
$cell.Style = "Normal"
$cell = $ws.Range("E24")
$cell.Value = "'  +1.78%  "
$cell.Style = "Normal"
$cell = $ws.Range("E25")
$cell.Value = "'  -0.27%  "
$cell.Style = "Normal"
$cell = $ws.Range("E26")
$cell.Value = "'  +1.01%  "
$cell.Style = "Normal"
$cell = $ws.Range("E27")
$cell.Value = "'  +0.28%  "
$cell.Style = "Normal"
$cell = $ws.Range("E28")
$cell.Value = "'  +2.65%  "
$cell.Style = "Normal"
$cell = $ws.Range("B29")
$cell.Value = "'USDe"
$cell.Style = "Normal"
$cell = $ws.Range("C29")
$cell.Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$cell.Style = "Normal"
$cell = $ws.Range("D29")
$cell.Value = "'0.999"
$cell.Style = "Normal"
$cell = $ws.Range("E29")
$cell.Value = "'  +0.17%  "
$cell.Style = "Normal"
$cell = $ws.Range("B30")
$cell.Value = "'PEPE"
$cell.Style = "Normal"
$cell = $ws.Range("C30")
$cell.Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$cell.Style = "Normal"
$cell = $ws.Range("D30")
$cell.Value = "'0.0₃0739"
$cell.Style = "Normal"
$cell = $ws.Range("E30")
$cell.Value = "'  +2.29%  "
$cell.Style = "Normal"
$cell = $ws.Range("E31")
$cell.Value = "'  +5.09%  "
$cell.Style = "Normal"
$cell = $ws.Range("E32")
$cell.Value = "'  -1.94%  "
$cell.Style = "Normal"
$cell = $ws.Range("D33")
$cell.Value = "'18.81"
$cell.Style = "Normal"
$cell = $ws.Range("E33")
$cell.Value = "'  +0.41%  "
$cell.Style = "Normal"
$cell = $ws.Range("D34")
$cell.Value = "'149.04"
$cell.Style = "Normal"
$cell = $ws.Range("E34")
$cell.Value = "'  -0.22%  "
$cell.Style = "Normal"
$cell = $ws.Range("D35")
$cell.Value = "'3.97"
$cell.Style = "Normal"
$cell = $ws.Range("E35")
$cell.Value = "'  -0.27%  "
$cell.Style = "Normal"
$cell = $ws.Range("E36")
$cell.Value = "'  -0.74%  "
$cell.Style = "Normal"
$cell = $ws.Range("D37")
$cell.Value = "'36.89"
$cell.Style = "Normal"
$cell = $ws.Range("E37")
$cell.Value = "'  +1.53%  "
$cell.Style = "Normal"
$cell = $ws.Range("E38")
$cell.Value = "'  +0.40%  "
$cell.Style = "Normal"
$cell = $ws.Range("B39")
$cell.Value = "'Fetch.AI"
$cell.Style = "Normal"
$cell = $ws.Range("C39")
$cell.Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$cell.Style = "Normal"
$cell = $ws.Range("E39")
$cell.Value = "'  +0.30%  "
$cell.Style = "Normal"
$cell = $ws.Range("B40")
$cell.Value = "'SuiNetwork"
$cell.Style = "Normal"
$cell = $ws.Range("C40")
$cell.Value = "'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$cell.Style = "Normal"
$cell = $ws.Range("D40")
$cell.Value = "'0.827"
$cell.Style = "Normal"
$cell = $ws.Range("E40")
$cell.Value = "'  -0.35%  "
$cell.Style = "Normal"
$cell = $ws.Range("E41")
$cell.Value = "'  +0.42%  "
$cell.Style = "Normal"
$cell = $ws.Range("E42")
$cell.Value = "'  +0.34%  "
$cell.Style = "Normal"
$cell = $ws.Range("D43")
$cell.Value = "'274.44"
$cell.Style = "Normal"
$cell = $ws.Range("E43")
$cell.Value = "'  -0.44%  "
$cell.Style = "Normal"
$cell = $ws.Range("D44")
$cell.Value = "'0.597"
$cell.Style = "Normal"
$cell = $ws.Range("E44")
$cell.Value = "'  +0.58%  "
$cell.Style = "Normal"
$cell = $ws.Range("E45")
$cell.Value = "'  +1.29%  "
$cell.Style = "Normal"
$cell = $ws.Range("E46")
$cell.Value = "'  +0.05%  "
$cell.Style = "Normal"
$cell = $ws.Range("D47")
$cell.Value = "'0.0523"
$cell.Style = "Normal"
$cell = $ws.Range("E47")
$cell.Value = "'  +0.63%  "
$cell.Style = "Normal"
$cell = $ws.Range("D48")
$cell.Value = "'1.947.32"
$cell.Style = "Normal"
$cell = $ws.Range("E48")
$cell.Value = "'  -1.92%  "
$cell.Style = "Normal"
$cell = $ws.Range("E49")
$cell.Value = "'  +1.52%  "
$cell.Style = "Normal"
$cell = $ws.Range("D50")
$cell.Value = "'18.28"
$cell.Style = "Normal"
$cell = $ws.Range("E50")
$cell.Value = "'  +1.30%  "
$cell.Style = "Normal"
$cell = $ws.Range("D51")
$cell.Value = "'4.50"
$cell.Style = "Normal"
$cell = $ws.Range("E51")
$cell.Value = "'  -1.96%  "
$cell.Style = "Normal"
